$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 0.907169
$ws.Cells.Item(2, 8).Value = 2.721507
$ws.Cells.Item(2, 9).Value = 0.001369063862079057
$ws.Cells.Item(2, 10).Value = 0.001369063862079057
$ws.Cells.Item(2, 11).Value = 1
$ws.Cells.Item(2, 12).Value = 0.3333333333333333
$ws.Cells.Item(2, 13).Value = 0.065175
$ws.Cells.Item(2, 14).Value = 0.195525
$ws.Cells.Item(2, 15).Value = 0.009404016458916581
$ws.Cells.Item(2, 16).Value = 0.009404016458916581
$ws.Cells.Item(2, 17).Value = 0.05912473957499999
$ws.Cells.Item(2, 18).Value = 0.532122656175
$ws.Cells.Item(2, 19).Value = 0.00001287469909229936
$ws.Cells.Item(2, 20).Value = 0.00001287469909229936

$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 0.907169
$ws.Cells.Item(3, 8).Value = 2.721507
$ws.Cells.Item(3, 9).Value = 0.001369063862079057
$ws.Cells.Item(3, 10).Value = 0.001369063862079057
$ws.Cells.Item(3, 13).Value = 6.718514333333332
$ws.Cells.Item(3, 14).Value = 20.155543
$ws.Cells.Item(3, 15).Value = 0.969405744075698
$ws.Cells.Item(3, 16).Value = 0.969405744075698
$ws.Cells.Item(3, 17).Value = 6.094827929255666
$ws.Cells.Item(3, 18).Value = 54.85345136330099
$ws.Cells.Item(3, 19).Value = 0.001327178371905897
$ws.Cells.Item(3, 20).Value = 0.001327178371905897

$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 0.907169
$ws.Cells.Item(4, 8).Value = 2.721507
$ws.Cells.Item(4, 9).Value = 0.001369063862079057
$ws.Cells.Item(4, 10).Value = 0.001369063862079057
$ws.Cells.Item(4, 11).Value = 1
$ws.Cells.Item(4, 12).Value = 0.3333333333333333
$ws.Cells.Item(4, 13).Value = 0.14686
$ws.Cells.Item(4, 14).Value = 0.44058
$ws.Cells.Item(4, 15).Value = 0.02119023946538534
$ws.Cells.Item(4, 16).Value = 0.02119023946538533
$ws.Cells.Item(4, 17).Value = 0.13322683934
$ws.Cells.Item(4, 18).Value = 1.19904155406
$ws.Cells.Item(4, 19).Value = 0.00002901079108086051
$ws.Cells.Item(4, 20).Value = 0.0000290107910808605

$ws.Cells.Item(5, 9).Value = 0.002841027838709403
$ws.Cells.Item(5, 10).Value = 0.002841027838709403
$ws.Cells.Item(5, 11).Value = 1
$ws.Cells.Item(5, 12).Value = 0.3333333333333333
$ws.Cells.Item(5, 13).Value = 0.065175
$ws.Cells.Item(5, 14).Value = 0.195525
$ws.Cells.Item(5, 15).Value = 0.009404016458916581
$ws.Cells.Item(5, 16).Value = 0.009404016458916581
$ws.Cells.Item(5, 17).Value = 0.122693349625
$ws.Cells.Item(5, 18).Value = 1.104240146625
$ws.Cells.Item(5, 19).Value = 0.00002671707255546342
$ws.Cells.Item(5, 20).Value = 0.00002671707255546343

$ws.Cells.Item(6, 9).Value = 0.002841027838709403
$ws.Cells.Item(6, 10).Value = 0.002841027838709403
$ws.Cells.Item(6, 13).Value = 6.718514333333332
$ws.Cells.Item(6, 14).Value = 20.155543
$ws.Cells.Item(6, 15).Value = 0.969405744075698
$ws.Cells.Item(6, 16).Value = 0.969405744075698
$ws.Cells.Item(6, 17).Value = 12.64774880031055
$ws.Cells.Item(6, 18).Value = 113.829739202795
$ws.Cells.Item(6, 19).Value = 0.002754108705923861
$ws.Cells.Item(6, 20).Value = 0.002754108705923861

$ws.Cells.Item(7, 9).Value = 0.002841027838709403
$ws.Cells.Item(7, 10).Value = 0.002841027838709403
$ws.Cells.Item(7, 11).Value = 1
$ws.Cells.Item(7, 12).Value = 0.3333333333333333
$ws.Cells.Item(7, 13).Value = 0.14686
$ws.Cells.Item(7, 14).Value = 0.44058
$ws.Cells.Item(7, 15).Value = 0.02119023946538534
$ws.Cells.Item(7, 16).Value = 0.02119023946538533
$ws.Cells.Item(7, 17).Value = 0.2764671319666667
$ws.Cells.Item(7, 18).Value = 2.4882041877
$ws.Cells.Item(7, 19).Value = 0.00006020206023007839
$ws.Cells.Item(7, 20).Value = 0.00006020206023007839

$ws.Cells.Item(8, 7).Value = 84.26343166666668
$ws.Cells.Item(8, 8).Value = 252.790295
$ws.Cells.Item(8, 9).Value = 0.1271670650006795
$ws.Cells.Item(8, 10).Value = 0.1271670650006795
$ws.Cells.Item(8, 11).Value = 1
$ws.Cells.Item(8, 12).Value = 0.3333333333333333
$ws.Cells.Item(8, 13).Value = 0.065175
$ws.Cells.Item(8, 14).Value = 0.195525
$ws.Cells.Item(8, 15).Value = 0.009404016458916581
$ws.Cells.Item(8, 16).Value = 0.009404016458916581
$ws.Cells.Item(8, 17).Value = 5.491869158875001
$ws.Cells.Item(8, 18).Value = 49.426822429875
$ws.Cells.Item(8, 19).Value = 0.001195881172298505
$ws.Cells.Item(8, 20).Value = 0.001195881172298505

$ws.Cells.Item(9, 7).Value = 84.26343166666668
$ws.Cells.Item(9, 8).Value = 252.790295
$ws.Cells.Item(9, 9).Value = 0.1271670650006795
$ws.Cells.Item(9, 10).Value = 0.1271670650006795
$ws.Cells.Item(9, 13).Value = 6.718514333333332
$ws.Cells.Item(9, 14).Value = 20.155543
$ws.Cells.Item(9, 15).Value = 0.969405744075698
$ws.Cells.Item(9, 16).Value = 0.969405744075698
$ws.Cells.Item(9, 17).Value = 566.1250734283539
$ws.Cells.Item(9, 18).Value = 5095.125660855185
$ws.Cells.Item(9, 19).Value = 0.1232764832689063
$ws.Cells.Item(9, 20).Value = 0.1232764832689063

$ws.Cells.Item(10, 7).Value = 84.26343166666668
$ws.Cells.Item(10, 8).Value = 252.790295
$ws.Cells.Item(10, 9).Value = 0.1271670650006795
$ws.Cells.Item(10, 10).Value = 0.1271670650006795
$ws.Cells.Item(10, 11).Value = 1
$ws.Cells.Item(10, 12).Value = 0.3333333333333333
$ws.Cells.Item(10, 13).Value = 0.14686
$ws.Cells.Item(10, 14).Value = 0.44058
$ws.Cells.Item(10, 15).Value = 0.02119023946538534
$ws.Cells.Item(10, 16).Value = 0.02119023946538533
$ws.Cells.Item(10, 17).Value = 12.37492757456667
$ws.Cells.Item(10, 18).Value = 111.3743481711
$ws.Cells.Item(10, 19).Value = 0.002694700559474621
$ws.Cells.Item(10, 20).Value = 0.00269470055947462

$ws.Cells.Item(11, 5).Value = 2
$ws.Cells.Item(11, 6).Value = 0.6666666666666666
$ws.Cells.Item(11, 7).Value = 0.244783
$ws.Cells.Item(11, 8).Value = 0.7343489999999999
$ws.Cells.Item(11, 9).Value = 0.0003694168995537743
$ws.Cells.Item(11, 10).Value = 0.0003694168995537743
$ws.Cells.Item(11, 11).Value = 1
$ws.Cells.Item(11, 12).Value = 0.3333333333333333
$ws.Cells.Item(11, 13).Value = 0.065175
$ws.Cells.Item(11, 14).Value = 0.195525
$ws.Cells.Item(11, 15).Value = 0.009404016458916581
$ws.Cells.Item(11, 16).Value = 0.009404016458916581
$ws.Cells.Item(11, 17).Value = 0.015953732025
$ws.Cells.Item(11, 18).Value = 0.143583588225
$ws.Cells.Item(11, 19).Value = 0.000003474002603605627
$ws.Cells.Item(11, 20).Value = 0.000003474002603605627

$ws.Cells.Item(12, 5).Value = 2
$ws.Cells.Item(12, 6).Value = 0.6666666666666666
$ws.Cells.Item(12, 7).Value = 0.244783
$ws.Cells.Item(12, 8).Value = 0.7343489999999999
$ws.Cells.Item(12, 9).Value = 0.0003694168995537743
$ws.Cells.Item(12, 10).Value = 0.0003694168995537743
$ws.Cells.Item(12, 13).Value = 6.718514333333332
$ws.Cells.Item(12, 14).Value = 20.155543
$ws.Cells.Item(12, 15).Value = 0.969405744075698
$ws.Cells.Item(12, 16).Value = 0.969405744075698
$ws.Cells.Item(12, 17).Value = 1.644578094056333
$ws.Cells.Item(12, 18).Value = 14.801202846507
$ws.Cells.Item(12, 19).Value = 0.000358114864386064
$ws.Cells.Item(12, 20).Value = 0.000358114864386064

$ws.Cells.Item(13, 5).Value = 2
$ws.Cells.Item(13, 6).Value = 0.6666666666666666
$ws.Cells.Item(13, 7).Value = 0.244783
$ws.Cells.Item(13, 8).Value = 0.7343489999999999
$ws.Cells.Item(13, 9).Value = 0.0003694168995537743
$ws.Cells.Item(13, 10).Value = 0.0003694168995537743
$ws.Cells.Item(13, 11).Value = 1
$ws.Cells.Item(13, 12).Value = 0.3333333333333333
$ws.Cells.Item(13, 13).Value = 0.14686
$ws.Cells.Item(13, 14).Value = 0.44058
$ws.Cells.Item(13, 15).Value = 0.02119023946538534
$ws.Cells.Item(13, 16).Value = 0.02119023946538533
$ws.Cells.Item(13, 17).Value = 0.03594883138
$ws.Cells.Item(13, 18).Value = 0.32353948242
$ws.Cells.Item(13, 19).Value = 0.000007828032564104679
$ws.Cells.Item(13, 20).Value = 0.000007828032564104677

$ws.Cells.Item(14, 7).Value = 62.79827133333333
$ws.Cells.Item(14, 8).Value = 188.394814
$ws.Cells.Item(14, 9).Value = 0.09477268720988248
$ws.Cells.Item(14, 10).Value = 0.09477268720988248
$ws.Cells.Item(14, 11).Value = 1
$ws.Cells.Item(14, 12).Value = 0.3333333333333333
$ws.Cells.Item(14, 13).Value = 0.065175
$ws.Cells.Item(14, 14).Value = 0.195525
$ws.Cells.Item(14, 15).Value = 0.009404016458916581
$ws.Cells.Item(14, 16).Value = 0.009404016458916581
$ws.Cells.Item(14, 17).Value = 4.09287733415
$ws.Cells.Item(14, 18).Value = 36.83589600735
$ws.Cells.Item(14, 19).Value = 0.0008912439103774878
$ws.Cells.Item(14, 20).Value = 0.0008912439103774878

$ws.Cells.Item(15, 7).Value = 62.79827133333333
$ws.Cells.Item(15, 8).Value = 188.394814
$ws.Cells.Item(15, 9).Value = 0.09477268720988248
$ws.Cells.Item(15, 10).Value = 0.09477268720988248
$ws.Cells.Item(15, 13).Value = 6.718514333333332
$ws.Cells.Item(15, 14).Value = 20.155543
$ws.Cells.Item(15, 15).Value = 0.969405744075698
$ws.Cells.Item(15, 16).Value = 0.969405744075698
$ws.Cells.Item(15, 17).Value = 421.9110860615557
$ws.Cells.Item(15, 18).Value = 3797.199774554002
$ws.Cells.Item(15, 19).Value = 0.09187318736274952
$ws.Cells.Item(15, 20).Value = 0.09187318736274952

$ws.Cells.Item(16, 7).Value = 62.79827133333333
$ws.Cells.Item(16, 8).Value = 188.394814
$ws.Cells.Item(16, 9).Value = 0.09477268720988248
$ws.Cells.Item(16, 10).Value = 0.09477268720988248
$ws.Cells.Item(16, 11).Value = 1
$ws.Cells.Item(16, 12).Value = 0.3333333333333333
$ws.Cells.Item(16, 13).Value = 0.14686
$ws.Cells.Item(16, 14).Value = 0.44058
$ws.Cells.Item(16, 15).Value = 0.02119023946538534
$ws.Cells.Item(16, 16).Value = 0.02119023946538533
$ws.Cells.Item(16, 17).Value = 9.222554128013334
$ws.Cells.Item(16, 18).Value = 83.00298715212
$ws.Cells.Item(16, 19).Value = 0.002008255936755472
$ws.Cells.Item(16, 20).Value = 0.002008255936755472

$ws.Cells.Item(17, 7).Value = 512.5237530000001
$ws.Cells.Item(17, 8).Value = 1537.571259
$ws.Cells.Item(17, 9).Value = 0.7734807391890958
$ws.Cells.Item(17, 10).Value = 0.7734807391890958
$ws.Cells.Item(17, 11).Value = 1
$ws.Cells.Item(17, 12).Value = 0.3333333333333333
$ws.Cells.Item(17, 13).Value = 0.065175
$ws.Cells.Item(17, 14).Value = 0.195525
$ws.Cells.Item(17, 15).Value = 0.009404016458916581
$ws.Cells.Item(17, 16).Value = 0.009404016458916581
$ws.Cells.Item(17, 17).Value = 33.403735601775
$ws.Cells.Item(17, 18).Value = 300.633620415975
$ws.Cells.Item(17, 19).Value = 0.00727382560198922
$ws.Cells.Item(17, 20).Value = 0.00727382560198922

$ws.Cells.Item(18, 7).Value = 512.5237530000001
$ws.Cells.Item(18, 8).Value = 1537.571259
$ws.Cells.Item(18, 9).Value = 0.7734807391890958
$ws.Cells.Item(18, 10).Value = 0.7734807391890958
$ws.Cells.Item(18, 13).Value = 6.718514333333332
$ws.Cells.Item(18, 14).Value = 20.155543
$ws.Cells.Item(18, 15).Value = 0.969405744075698
$ws.Cells.Item(18, 16).Value = 0.969405744075698
$ws.Cells.Item(18, 17).Value = 3443.398180704293
$ws.Cells.Item(18, 18).Value = 30990.58362633863
$ws.Cells.Item(18, 19).Value = 0.7498166715018263
$ws.Cells.Item(18, 20).Value = 0.7498166715018263

$ws.Cells.Item(19, 7).Value = 512.5237530000001
$ws.Cells.Item(19, 8).Value = 1537.571259
$ws.Cells.Item(19, 9).Value = 0.7734807391890958
$ws.Cells.Item(19, 10).Value = 0.7734807391890958
$ws.Cells.Item(19, 11).Value = 1
$ws.Cells.Item(19, 12).Value = 0.3333333333333333
$ws.Cells.Item(19, 13).Value = 0.14686
$ws.Cells.Item(19, 14).Value = 0.44058
$ws.Cells.Item(19, 15).Value = 0.02119023946538534
$ws.Cells.Item(19, 16).Value = 0.02119023946538533
$ws.Cells.Item(19, 17).Value = 75.26923836558002
$ws.Cells.Item(19, 18).Value = 677.4231452902201
$ws.Cells.Item(19, 19).Value = 0.0163902420852802
$ws.Cells.Item(19, 20).Value = 0.0163902420852802
